# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Murcott, $/caja 15 kilos granel) at rows 102-103,
# pushing all subsequent rows down by two (dimension grows from T132 to T134).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("102:103").Insert()

# Row 102
$ws.Cells.Item(102,1).Value2 = 11
$ws.Cells.Item(102,2).Value2 = 'Vega Monumental Concepción'
$ws.Cells.Item(102,3).Value2 = 'Bíobío'
$ws.Cells.Item(102,4).Value2 = 44782
$ws.Cells.Item(102,5).Value2 = 8
$ws.Cells.Item(102,6).Value2 = 'Fruta'
$ws.Cells.Item(102,7).Value2 = 100102
$ws.Cells.Item(102,8).Value2 = 'Cítricos'
$ws.Cells.Item(102,9).Value2 = 100102004
$ws.Cells.Item(102,10).Value2 = 'Mandarina'
$ws.Cells.Item(102,11).Value2 = 'Murcott'
$ws.Cells.Item(102,12).Value2 = 'Primera'
$ws.Cells.Item(102,13).Value2 = 250
$ws.Cells.Item(102,14).Value2 = 100
$ws.Cells.Item(102,15).Value2 = 120
$ws.Cells.Item(102,16).Value2 = 108
$ws.Cells.Item(102,17).Value2 = '$/caja 15 kilos granel'
$ws.Cells.Item(102,18).Value2 = 'Región de O''Higgins'
$ws.Cells.Item(102,19).Value2 = 7
$ws.Cells.Item(102,20).Value2 = 15

# Row 103
$ws.Cells.Item(103,1).Value2 = 11
$ws.Cells.Item(103,2).Value2 = 'Vega Monumental Concepción'
$ws.Cells.Item(103,3).Value2 = 'Bíobío'
$ws.Cells.Item(103,4).Value2 = 44782
$ws.Cells.Item(103,5).Value2 = 8
$ws.Cells.Item(103,6).Value2 = 'Fruta'
$ws.Cells.Item(103,7).Value2 = 100102
$ws.Cells.Item(103,8).Value2 = 'Cítricos'
$ws.Cells.Item(103,9).Value2 = 100102004
$ws.Cells.Item(103,10).Value2 = 'Mandarina'
$ws.Cells.Item(103,11).Value2 = 'Murcott'
$ws.Cells.Item(103,12).Value2 = 'Segunda'
$ws.Cells.Item(103,13).Value2 = 200
$ws.Cells.Item(103,14).Value2 = 7000
$ws.Cells.Item(103,15).Value2 = 7500
$ws.Cells.Item(103,16).Value2 = 7250
$ws.Cells.Item(103,17).Value2 = '$/caja 15 kilos granel'
$ws.Cells.Item(103,18).Value2 = 'Región de O''Higgins'
$ws.Cells.Item(103,19).Value2 = 483
$ws.Cells.Item(103,20).Value2 = 15

